$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace "Cruzer Blade / SanDisk" USB drive spec with "Vostro / Dell" machine spec ---
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"

# --- Row 3: Arabic translation row, same new data ---
$ws.Range("B3").Value = "ستر  "
$ws.Range("C3").Value = "دلّ  "
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("G3").Value = "لأخذ التسجيلات"

# --- View state: scroll sheet so column C is leftmost, and select from row 4 down ---
$ws.Range("C4").Select() | Out-Null
$ws.Rows("4:1048576").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

# --- Page setup: A4 paper, portrait orientation ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
